$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4236.8335
$ws.Range("J17").Value = 4855.125
$ws.Range("L17").Value = 14565.375
$ws.Range("N17").Value = -14901.375
$ws.Range("H76").Value = 3099.8572
$ws.Range("I76").Value = 2866.5
$ws.Range("K76").Value = 2866.5
$ws.Range("M76").Value = -2551.5
$ws.Range("H79").Value = 3099.8572
$ws.Range("I79").Value = 2866.5
$ws.Range("K79").Value = 2866.5
$ws.Range("M79").Value = -1774.5
$ws.Range("H107").Value = 515.8889
$ws.Range("I107").Value = 515.8889
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 515.8889
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1404.1111
$ws.Range("N107").ClearContents()
$ws.Range("H116").Value = 4417.4546
$ws.Range("I116").Value = 3707.375
$ws.Range("J116").Value = 6311
$ws.Range("K116").Value = 3707.375
$ws.Range("L116").Value = 6311
$ws.Range("M116").Value = -265.375
$ws.Range("N116").Value = -13195
$ws.Range("H118").Value = 947.63635
$ws.Range("I118").Value = 936.55554
$ws.Range("J118").Value = 997.5
$ws.Range("K118").Value = 2809.66662
$ws.Range("L118").Value = 2992.5
$ws.Range("M118").Value = -1152.66662
$ws.Range("N118").Value = -6306.5
$ws.Range("H125").Value = 8349.75
$ws.Range("I125").Value = 4749.75
$ws.Range("K125").Value = 42747.75
$ws.Range("M125").Value = -40287.75
$ws.Range("H129").Value = 2597.5
$ws.Range("I129").Value = 1066.6
$ws.Range("K129").Value = 3199.8
$ws.Range("M129").Value = 1800.2

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 24000
$ws.Range("I53").Value = 24000
$ws.Range("K53").Value = 24000
$ws.Range("M53").Value = -23318
$ws.Range("H132").Value = 2838.7273
$ws.Range("I132").Value = 2838.7273
$ws.Range("K132").Value = 8516.1819
$ws.Range("M132").Value = -5986.1819

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1580.9642
$ws.Range("J86").Value = 3268.1428
$ws.Range("L86").Value = 3268.1428
$ws.Range("N86").Value = -5514.1428
$ws.Range("H89").Value = 1580.9642
$ws.Range("J89").Value = 3268.1428
$ws.Range("L89").Value = 16340.714
$ws.Range("N89").Value = -27572.714
$ws.Range("H107").Value = 1284.1428
$ws.Range("I107").Value = 1284.1428
$ws.Range("K107").Value = 1284.1428
$ws.Range("M107").Value = 635.8571999999999
$ws.Range("H134").Value = 9937.77
$ws.Range("I134").Value = 11021.333
$ws.Range("K134").Value = 33063.999
$ws.Range("M134").Value = -30528.999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2486.25
$ws.Range("I16").Value = 1190.8
$ws.Range("K16").Value = 1190.8
$ws.Range("M16").Value = -903.8
$ws.Range("H22").Value = 261.85715
$ws.Range("I22").Value = 232.75
$ws.Range("K22").Value = 232.75
$ws.Range("M22").Value = 117.25
$ws.Range("H41").Value = 15816.333
$ws.Range("J41").Value = 23333.334
$ws.Range("L41").Value = 23333.334
$ws.Range("N41").Value = -24189.334
$ws.Range("H94").Value = 1706
$ws.Range("J94").Value = 1161.5
$ws.Range("L94").Value = 1161.5
$ws.Range("N94").Value = -2063.5
$ws.Range("H105").Value = 763.4737
$ws.Range("I105").Value = 717.8125
$ws.Range("K105").Value = 717.8125
$ws.Range("M105").Value = 1029.1875
$ws.Range("H113").Value = 2486.25
$ws.Range("I113").Value = 1190.8
$ws.Range("K113").Value = 1190.8
$ws.Range("M113").Value = 979.2
$ws.Range("H134").Value = 5872.75
$ws.Range("I134").Value = 2830.3333
$ws.Range("K134").Value = 8490.999899999999
$ws.Range("M134").Value = -5955.999899999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 362.83334
$ws.Range("J23").Value = 676.6667
$ws.Range("L23").Value = 2030.0001
$ws.Range("N23").Value = -2500.0001
$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("J25").Value = 1000
$ws.Range("K25").Value = 3000
$ws.Range("L25").Value = 3000
$ws.Range("M25").Value = -2831
$ws.Range("N25").Value = -3338
$ws.Range("H30").Value = 1000
$ws.Range("I30").Value = 1000
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 3000
$ws.Range("L30").Value = 3000
$ws.Range("M30").Value = -2898
$ws.Range("N30").Value = -3204
$ws.Range("H75").Value = 1265
$ws.Range("I75").Value = 30
$ws.Range("J75").Value = 2500
$ws.Range("K75").Value = 90
$ws.Range("L75").Value = 7500
$ws.Range("M75").Value = 908
$ws.Range("N75").Value = -9496
$ws.Range("H78").Value = 1265
$ws.Range("I78").Value = 30
$ws.Range("J78").Value = 2500
$ws.Range("K78").Value = 270
$ws.Range("L78").Value = 22500
$ws.Range("M78").Value = 4722
$ws.Range("N78").Value = -32484
$ws.Range("H103").Value = 725
$ws.Range("I103").Value = 725
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 2175
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -1296
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 632.125
$ws.Range("I113").Value = 522
$ws.Range("K113").Value = 1566
$ws.Range("M113").Value = 604
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("M114").ClearContents()
$ws.Range("N114").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2691.0967
$ws.Range("I80").Value = 1820.8334
$ws.Range("J80").Value = 2899.96
$ws.Range("K80").Value = 1820.8334
$ws.Range("L80").Value = 2899.96
$ws.Range("M80").Value = -822.8334
$ws.Range("N80").Value = -4895.96
$ws.Range("H83").Value = 2691.0967
$ws.Range("I83").Value = 1820.8334
$ws.Range("J83").Value = 2899.96
$ws.Range("K83").Value = 9104.166999999999
$ws.Range("L83").Value = 14499.8
$ws.Range("M83").Value = -4112.166999999999
$ws.Range("N83").Value = -24483.8
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1078.4
$ws.Range("I31").Value = 596
$ws.Range("K31").Value = 596
$ws.Range("M31").Value = -348
$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 31710.75
$ws.Range("I61").Value = 31710.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 31710.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -31418.75
$ws.Range("N61").ClearContents()
$ws.Range("H62").Value = 2297.5
$ws.Range("I62").Value = 2277.2
$ws.Range("K62").Value = 2277.2
$ws.Range("M62").Value = -1653.2
$ws.Range("H65").Value = 2297.5
$ws.Range("I65").Value = 2277.2
$ws.Range("K65").Value = 11386
$ws.Range("M65").Value = -8266
$ws.Range("H126").Value = 2869.8
$ws.Range("I126").Value = 2440.6365
$ws.Range("J126").Value = 4050
$ws.Range("K126").Value = 7321.9095
$ws.Range("L126").Value = 12150
$ws.Range("M126").Value = -4851.9095
$ws.Range("N126").Value = -17090
$ws.Range("H132").Value = 1216.6666
$ws.Range("I132").Value = 1216.6666
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3649.9998
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1119.9998
$ws.Range("N132").ClearContents()
